$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 10 (student "Губанов Арсений") homework grades.
# Column F already had the highlighted (yellow/green) fill used for "missing"
# entries; match the plain bordered style used by columns G/H/I before filling
# in the grade values.
$ws.Range("G10").Copy($ws.Range("F10"))

$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 5
$ws.Range("H10").Value = 5
$ws.Range("I10").Value = 5

# Move the active selection to I10 (as recorded in the saved view state).
$ws.Range("I10").Select()
